# "finish the function of bless system" (close #97)
#
# The Scene sheet's QuestRandom column (F) gets a new "witchhome;N" random
# reward appended to four rows' pipe-delimited reward lists. (All other
# shared-string index churn visible in the raw OOXML diff is just the
# sharedStrings table being re-packed around these new strings - the
# underlying cell text for every other cell is unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")

$ws.Range("F8").Value  = "met;30|treasure;25|witchhome;20"
$ws.Range("F16").Value = "forestfire;35|witchhome;10"
$ws.Range("F18").Value = "met;30|witchhome;30"
$ws.Range("F22").Value = "met;30|forestfire;20|witchhome;40"

# Match the author's final cell selection in the saved workbook.
$ws.Range("F16").Select()
